$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set A1 to the new header text
$ws.Range("A1").Value = "Вопросы"

# Remove the now-unneeded rows 2-4 entirely (also clears their content/formatting)
$ws.Rows("2:4").Delete()
